$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.684.55'
$ws.Range("E2").Value = '  -4.77%  '

$ws.Range("D3").Value = '2.313.34'
$ws.Range("E3").Value = '  -6.25%  '

$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").Value = '''305.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.65%  '

$ws.Range("D6").Value = '''83.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.59%  '

$ws.Range("E7").Value = '  -3.83%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -5.32%  '

$ws.Range("D10").Value = '''0.0806'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.56%  '

$ws.Range("D11").Value = '''29.74'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.31%  '

$ws.Range("D12").Value = '''0.109'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").Value = '2.669.59'
$ws.Range("E13").Value = '  -6.25%  '

$ws.Range("E14").Value = '  -7.34%  '

$ws.Range("D15").Value = '''14.58'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.73%  '

$ws.Range("D16").Value = '2.318.85'
$ws.Range("E16").Value = '  -5.95%  '

$ws.Range("D17").Value = '''0.748'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.18%  '

$ws.Range("D18").Value = '39.680.18'
$ws.Range("E18").Value = '  -4.67%  '

$ws.Range("D19").Value = '0.0₃0895'
$ws.Range("E19").Value = '  -4.68%  '

$ws.Range("E20").Value = '  -6.20%  '

$ws.Range("D21").Value = '''67.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.75%  '

$ws.Range("E22").Value = '  -6.17%  '

$ws.Range("D23").Value = '''233.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.21%  '

$ws.Range("D24").Value = '''2.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.11%  '

$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").Value = '''1.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.61%  '

$ws.Range("D27").Value = '''22.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.42%  '

$ws.Range("E28").Value = '  -2.11%  '

$ws.Range("E29").Value = '  -5.17%  '

$ws.Range("D30").Value = '''34.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.47%  '

$ws.Range("D31").Value = '''150.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.20%  '

$ws.Range("E32").Value = '  -0.19%  '

$ws.Range("E33").Value = '  -6.95%  '

$ws.Range("D34").Value = '''2.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.95%  '

$ws.Range("D35").Value = '''0.0715'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.30%  '

$ws.Range("E36").Value = '  -2.84%  '

$ws.Range("D37").Value = '''0.0985'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.28%  '

$ws.Range("D38").Value = '''2.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.42%  '

$ws.Range("D39").Value = '''15.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.22%  '

$ws.Range("E40").Value = '  -7.46%  '

$ws.Range("D41").Value = '''3.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.69%  '

$ws.Range("E42").Value = '  -2.56%  '

$ws.Range("D43").Value = '1.937.10'
$ws.Range("E43").Value = '  -3.07%  '

$ws.Range("D44").Value = '''0.0263'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.48%  '

$ws.Range("D45").Value = '''17.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.13%  '

$ws.Range("E46").Value = '  -1.95%  '

$ws.Range("E47").Value = '  -10.16%  '

$ws.Range("D48").Value = '2.542.61'
$ws.Range("E48").Value = '  -6.74%  '

$ws.Range("D49").Value = '''92.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.99%  '

$ws.Range("D50").Value = '''69.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.51%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '''62.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.24%  '
